$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# MOSIP-15883: add new master-template rows for UIN duplicate/technical-issue
# re-registration subject templates and UIN activation/deactivation/generation/update subjects.

# Row 438
$ws.Cells.Item(438,1).Value = 1274
$ws.Cells.Item(438,2).Value = "Registration Failed because you have already Registered"
$ws.Cells.Item(438,3).Value = "Registration Failed because you have already Registered"
$ws.Cells.Item(438,4).Value = "txt"
$ws.Cells.Item(438,5).Value = "velocity"
$ws.Cells.Item(438,6).Value = "Registration Failed because you have already Registered"
$ws.Cells.Item(438,6).WrapText = $true
$ws.Cells.Item(438,7).Value = 10003
$ws.Cells.Item(438,8).Value = "Registration Processor"
$ws.Cells.Item(438,9).Value = "RPR_DUP_UIN_EMAIL_SUB"
$ws.Cells.Item(438,10).Value = "eng"
$ws.Cells.Item(438,11).Value = $true
$ws.Cells.Item(438,12).Value = "superadmin"
$ws.Cells.Item(438,13).Value = "now()"

# Row 439
$ws.Cells.Item(439,1).Value = 1274
$ws.Cells.Item(439,2).Value = "L'inscription a échoué car vous êtes déjà inscrit"
$ws.Cells.Item(439,3).Value = "L'inscription a échoué car vous êtes déjà inscrit"
$ws.Cells.Item(439,4).Value = "txt"
$ws.Cells.Item(439,5).Value = "velocity"
$ws.Cells.Item(439,6).Value = "L'inscription a échoué car vous êtes déjà inscrit"
$ws.Cells.Item(439,6).WrapText = $true
$ws.Cells.Item(439,7).Value = 10003
$ws.Cells.Item(439,8).Value = "Registration Processor"
$ws.Cells.Item(439,9).Value = "RPR_DUP_UIN_EMAIL_SUB"
$ws.Cells.Item(439,10).Value = "fra"
$ws.Cells.Item(439,11).Value = $true
$ws.Cells.Item(439,12).Value = "superadmin"
$ws.Cells.Item(439,13).Value = "now()"

# Row 440
$ws.Cells.Item(440,1).Value = 1274
$ws.Cells.Item(440,2).Value = "فشل التسجيل لأنك قمت بالتسجيل بالفعل"
$ws.Cells.Item(440,3).Value = "فشل التسجيل لأنك قمت بالتسجيل بالفعل"
$ws.Cells.Item(440,4).Value = "txt"
$ws.Cells.Item(440,5).Value = "velocity"
$ws.Cells.Item(440,6).Value = "فشل التسجيل لأنك قمت بالتسجيل بالفعل"
$ws.Cells.Item(440,6).WrapText = $true
$ws.Cells.Item(440,7).Value = 10003
$ws.Cells.Item(440,8).Value = "Registration Processor"
$ws.Cells.Item(440,9).Value = "RPR_DUP_UIN_EMAIL_SUB"
$ws.Cells.Item(440,10).Value = "ara"
$ws.Cells.Item(440,11).Value = $true
$ws.Cells.Item(440,12).Value = "superadmin"
$ws.Cells.Item(440,13).Value = "now()"

# Row 441
$ws.Cells.Item(441,1).Value = 1275
$ws.Cells.Item(441,2).Value = "Re-Register because there was a Technical Issue"
$ws.Cells.Item(441,3).Value = "Re-Register because there was a Technical Issue"
$ws.Cells.Item(441,4).Value = "txt"
$ws.Cells.Item(441,5).Value = "velocity"
$ws.Cells.Item(441,6).Value = "Re-Register because there was a Technical Issue"
$ws.Cells.Item(441,6).WrapText = $true
$ws.Cells.Item(441,7).Value = 10003
$ws.Cells.Item(441,8).Value = "Registration Processor"
$ws.Cells.Item(441,9).Value = "RPR_TEC_ISSUE_EMAIL_SUB"
$ws.Cells.Item(441,10).Value = "eng"
$ws.Cells.Item(441,11).Value = $true
$ws.Cells.Item(441,12).Value = "superadmin"
$ws.Cells.Item(441,13).Value = "now()"

# Row 442
$ws.Cells.Item(442,1).Value = 1275
$ws.Cells.Item(442,2).Value = "Réinscrivez-vous car il y a eu un problème technique"
$ws.Cells.Item(442,3).Value = "Réinscrivez-vous car il y a eu un problème technique"
$ws.Cells.Item(442,4).Value = "txt"
$ws.Cells.Item(442,5).Value = "velocity"
$ws.Cells.Item(442,6).Value = "Réinscrivez-vous car il y a eu un problème technique"
$ws.Cells.Item(442,6).WrapText = $true
$ws.Cells.Item(442,7).Value = 10003
$ws.Cells.Item(442,8).Value = "Registration Processor"
$ws.Cells.Item(442,9).Value = "RPR_TEC_ISSUE_EMAIL_SUB"
$ws.Cells.Item(442,10).Value = "fra"
$ws.Cells.Item(442,11).Value = $true
$ws.Cells.Item(442,12).Value = "superadmin"
$ws.Cells.Item(442,13).Value = "now()"

# Row 443
$ws.Cells.Item(443,1).Value = 1275
$ws.Cells.Item(443,2).Value = "إعادة التسجيل بسبب وجود مشكلة فنية"
$ws.Cells.Item(443,3).Value = "إعادة التسجيل بسبب وجود مشكلة فنية"
$ws.Cells.Item(443,4).Value = "txt"
$ws.Cells.Item(443,5).Value = "velocity"
$ws.Cells.Item(443,6).Value = "إعادة التسجيل بسبب وجود مشكلة فنية"
$ws.Cells.Item(443,6).WrapText = $true
$ws.Cells.Item(443,7).Value = 10003
$ws.Cells.Item(443,8).Value = "Registration Processor"
$ws.Cells.Item(443,9).Value = "RPR_TEC_ISSUE_EMAIL_SUB"
$ws.Cells.Item(443,10).Value = "ara"
$ws.Cells.Item(443,11).Value = $true
$ws.Cells.Item(443,12).Value = "superadmin"
$ws.Cells.Item(443,13).Value = "now()"

# Row 444
$ws.Cells.Item(444,1).Value = 1276
$ws.Cells.Item(444,2).Value = "Uin is activated successfully"
$ws.Cells.Item(444,3).Value = "Uin is activated successfully"
$ws.Cells.Item(444,4).Value = "txt"
$ws.Cells.Item(444,5).Value = "velocity"
$ws.Cells.Item(444,6).Value = "Uin is activated successfully"
$ws.Cells.Item(444,7).Value = 10003
$ws.Cells.Item(444,8).Value = "Registration Processor"
$ws.Cells.Item(444,9).Value = "RPR_UIN_REAC_EMAIL_SUB"
$ws.Cells.Item(444,10).Value = "eng"
$ws.Cells.Item(444,11).Value = $true
$ws.Cells.Item(444,12).Value = "superadmin"
$ws.Cells.Item(444,13).Value = "now()"

# Row 445
$ws.Cells.Item(445,1).Value = 1276
$ws.Cells.Item(445,2).Value = "Uin est activé avec succès"
$ws.Cells.Item(445,3).Value = "Uin est activé avec succès"
$ws.Cells.Item(445,4).Value = "txt"
$ws.Cells.Item(445,5).Value = "velocity"
$ws.Cells.Item(445,6).Value = "Uin est activé avec succès"
$ws.Cells.Item(445,6).WrapText = $true
$ws.Cells.Item(445,7).Value = 10003
$ws.Cells.Item(445,8).Value = "Registration Processor"
$ws.Cells.Item(445,9).Value = "RPR_UIN_REAC_EMAIL_SUB"
$ws.Cells.Item(445,10).Value = "fra"
$ws.Cells.Item(445,11).Value = $true
$ws.Cells.Item(445,12).Value = "superadmin"
$ws.Cells.Item(445,13).Value = "now()"

# Row 446
$ws.Cells.Item(446,1).Value = 1276
$ws.Cells.Item(446,2).Value = "تم تفعيل Uin بنجاح"
$ws.Cells.Item(446,3).Value = "تم تفعيل Uin بنجاح"
$ws.Cells.Item(446,4).Value = "txt"
$ws.Cells.Item(446,5).Value = "velocity"
$ws.Cells.Item(446,6).Value = "تم تفعيل Uin بنجاح"
$ws.Cells.Item(446,6).WrapText = $true
$ws.Cells.Item(446,7).Value = 10003
$ws.Cells.Item(446,8).Value = "Registration Processor"
$ws.Cells.Item(446,9).Value = "RPR_UIN_REAC_EMAIL_SUB"
$ws.Cells.Item(446,10).Value = "ara"
$ws.Cells.Item(446,11).Value = $true
$ws.Cells.Item(446,12).Value = "superadmin"
$ws.Cells.Item(446,13).Value = "now()"

# Row 447
$ws.Cells.Item(447,1).Value = 1277
$ws.Cells.Item(447,2).Value = "Uin is deactivated"
$ws.Cells.Item(447,3).Value = "Uin is deactivated"
$ws.Cells.Item(447,4).Value = "txt"
$ws.Cells.Item(447,5).Value = "velocity"
$ws.Cells.Item(447,6).Value = "Uin is deactivated"
$ws.Cells.Item(447,6).WrapText = $true
$ws.Cells.Item(447,7).Value = 10003
$ws.Cells.Item(447,8).Value = "Registration Processor"
$ws.Cells.Item(447,9).Value = "RPR_UIN_DEAC_EMAIL_SUB"
$ws.Cells.Item(447,10).Value = "eng"
$ws.Cells.Item(447,11).Value = $true
$ws.Cells.Item(447,12).Value = "superadmin"
$ws.Cells.Item(447,13).Value = "now()"

# Row 448
$ws.Cells.Item(448,1).Value = 1277
$ws.Cells.Item(448,2).Value = "Uin est désactivé"
$ws.Cells.Item(448,3).Value = "Uin est désactivé"
$ws.Cells.Item(448,4).Value = "txt"
$ws.Cells.Item(448,5).Value = "velocity"
$ws.Cells.Item(448,6).Value = "Uin est désactivé"
$ws.Cells.Item(448,6).WrapText = $true
$ws.Cells.Item(448,7).Value = 10003
$ws.Cells.Item(448,8).Value = "Registration Processor"
$ws.Cells.Item(448,9).Value = "RPR_UIN_DEAC_EMAIL_SUB"
$ws.Cells.Item(448,10).Value = "fra"
$ws.Cells.Item(448,11).Value = $true
$ws.Cells.Item(448,12).Value = "superadmin"
$ws.Cells.Item(448,13).Value = "now()"

# Row 449
$ws.Cells.Item(449,1).Value = 1277
$ws.Cells.Item(449,2).Value = "تم إلغاء تنشيط Uin"
$ws.Cells.Item(449,3).Value = "تم إلغاء تنشيط Uin"
$ws.Cells.Item(449,4).Value = "txt"
$ws.Cells.Item(449,5).Value = "velocity"
$ws.Cells.Item(449,6).Value = "تم إلغاء تنشيط Uin"
$ws.Cells.Item(449,6).WrapText = $true
$ws.Cells.Item(449,7).Value = 10003
$ws.Cells.Item(449,8).Value = "Registration Processor"
$ws.Cells.Item(449,9).Value = "RPR_UIN_DEAC_EMAIL_SUB"
$ws.Cells.Item(449,10).Value = "ara"
$ws.Cells.Item(449,11).Value = $true
$ws.Cells.Item(449,12).Value = "superadmin"
$ws.Cells.Item(449,13).Value = "now()"

# Row 450
$ws.Cells.Item(450,1).Value = 1278
$ws.Cells.Item(450,2).Value = "UIN Generated"
$ws.Cells.Item(450,3).Value = "UIN Generated"
$ws.Cells.Item(450,4).Value = "txt"
$ws.Cells.Item(450,5).Value = "velocity"
$ws.Cells.Item(450,6).Value = "UIN Generated"
$ws.Cells.Item(450,6).WrapText = $true
$ws.Cells.Item(450,7).Value = 10003
$ws.Cells.Item(450,8).Value = "Registration Processor"
$ws.Cells.Item(450,9).Value = "RPR_UIN_GEN_EMAIL_SUB"
$ws.Cells.Item(450,10).Value = "eng"
$ws.Cells.Item(450,11).Value = $true
$ws.Cells.Item(450,12).Value = "superadmin"
$ws.Cells.Item(450,13).Value = "now()"

# Row 451
$ws.Cells.Item(451,1).Value = 1278
$ws.Cells.Item(451,2).Value = "UIN généré"
$ws.Cells.Item(451,3).Value = "UIN généré"
$ws.Cells.Item(451,4).Value = "txt"
$ws.Cells.Item(451,5).Value = "velocity"
$ws.Cells.Item(451,6).Value = "UIN généré"
$ws.Cells.Item(451,6).WrapText = $true
$ws.Cells.Item(451,7).Value = 10003
$ws.Cells.Item(451,8).Value = "Registration Processor"
$ws.Cells.Item(451,9).Value = "RPR_UIN_GEN_EMAIL_SUB"
$ws.Cells.Item(451,10).Value = "fra"
$ws.Cells.Item(451,11).Value = $true
$ws.Cells.Item(451,12).Value = "superadmin"
$ws.Cells.Item(451,13).Value = "now()"

# Row 452
$ws.Cells.Item(452,1).Value = 1278
$ws.Cells.Item(452,2).Value = "تم إنشاء UIN"
$ws.Cells.Item(452,3).Value = "تم إنشاء UIN"
$ws.Cells.Item(452,4).Value = "txt"
$ws.Cells.Item(452,5).Value = "velocity"
$ws.Cells.Item(452,6).Value = "تم إنشاء UIN"
$ws.Cells.Item(452,6).WrapText = $true
$ws.Cells.Item(452,7).Value = 10003
$ws.Cells.Item(452,8).Value = "Registration Processor"
$ws.Cells.Item(452,9).Value = "RPR_UIN_GEN_EMAIL_SUB"
$ws.Cells.Item(452,10).Value = "ara"
$ws.Cells.Item(452,11).Value = $true
$ws.Cells.Item(452,12).Value = "superadmin"
$ws.Cells.Item(452,13).Value = "now()"

# Row 453
$ws.Cells.Item(453,1).Value = 1279
$ws.Cells.Item(453,2).Value = "UIN Details Updated"
$ws.Cells.Item(453,3).Value = "UIN Details Updated"
$ws.Cells.Item(453,4).Value = "txt"
$ws.Cells.Item(453,5).Value = "velocity"
$ws.Cells.Item(453,6).Value = "UIN Details Updated"
$ws.Cells.Item(453,6).WrapText = $true
$ws.Cells.Item(453,7).Value = 10003
$ws.Cells.Item(453,8).Value = "Registration Processor"
$ws.Cells.Item(453,9).Value = "RPR_UIN_UPD_EMAIL_SUB"
$ws.Cells.Item(453,10).Value = "eng"
$ws.Cells.Item(453,11).Value = $true
$ws.Cells.Item(453,12).Value = "superadmin"
$ws.Cells.Item(453,13).Value = "now()"

# Row 454
$ws.Cells.Item(454,1).Value = 1279
$ws.Cells.Item(454,2).Value = "Détails UIN mis à jour"
$ws.Cells.Item(454,3).Value = "Détails UIN mis à jour"
$ws.Cells.Item(454,4).Value = "txt"
$ws.Cells.Item(454,5).Value = "velocity"
$ws.Cells.Item(454,6).Value = "Détails UIN mis à jour"
$ws.Cells.Item(454,6).WrapText = $true
$ws.Cells.Item(454,7).Value = 10003
$ws.Cells.Item(454,8).Value = "Registration Processor"
$ws.Cells.Item(454,9).Value = "RPR_UIN_UPD_EMAIL_SUB"
$ws.Cells.Item(454,10).Value = "fra"
$ws.Cells.Item(454,11).Value = $true
$ws.Cells.Item(454,12).Value = "superadmin"
$ws.Cells.Item(454,13).Value = "now()"

# Row 455
$ws.Cells.Item(455,1).Value = 1279
$ws.Cells.Item(455,2).Value = "تم تحديث تفاصيل UIN"
$ws.Cells.Item(455,3).Value = "تم تحديث تفاصيل UIN"
$ws.Cells.Item(455,4).Value = "txt"
$ws.Cells.Item(455,5).Value = "velocity"
$ws.Cells.Item(455,6).Value = "تم تحديث تفاصيل UIN"
$ws.Cells.Item(455,6).WrapText = $true
$ws.Cells.Item(455,7).Value = 10003
$ws.Cells.Item(455,8).Value = "Registration Processor"
$ws.Cells.Item(455,9).Value = "RPR_UIN_UPD_EMAIL_SUB"
$ws.Cells.Item(455,10).Value = "ara"
$ws.Cells.Item(455,11).Value = $true
$ws.Cells.Item(455,12).Value = "superadmin"
$ws.Cells.Item(455,13).Value = "now()"

[void]$ws.Range("F456").Select()
$excel.ActiveWindow.ScrollRow = 438
$excel.ActiveWindow.ScrollColumn = 1
